# Logged Week 15 and simulated Week 16
# Update row 3 ("R" - likely Road/Away totals) target depth data on both
# the OFF and DEF sheets with the latest cumulative stats.

$wb = $excel.ActiveWorkbook

# OFF sheet: Short Att, Short Comp, Deep Att, Deep Comp columns (B:E), row 3
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 198
$wsOff.Range("C3").Value = 141
$wsOff.Range("D3").Value = 57
$wsOff.Range("E3").Value = 27

# DEF sheet: Short Att, Short Comp, Deep Att, Deep Comp columns (B:E), row 3
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 174
$wsDef.Range("C3").Value = 127
$wsDef.Range("D3").Value = 49
$wsDef.Range("E3").Value = 19
